$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells G1/H1, styled like the existing header row (copy style from F1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Updated metric values in row 2
$ws.Range("B2").Value = 0.3442298297147584
$ws.Range("C2").Value = 0.9933777835200235
$ws.Range("D2").Value = 0.4656029464458803

# Updated model description text (multi-line)
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 AdaBoostRegressor(learning_rate=0.1, n_estimators=100))])"

# New Elapsed Time / CPU values
$ws.Range("G2").Value = 0.125854933266722
$ws.Range("H2").Value = 0.9890000000000001
